$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 19 (chunk 0)
$ws.Range("H19").Value = 2024.75
$ws.Range("I19").Value = 1966.3334
$ws.Range("K19").Value = 1966.3334
$ws.Range("M19").Value = -1791.3334
# row 28 (chunk 1)
$ws.Range("H28").Value = 1843.9131
$ws.Range("I28").Value = 1583.2354
$ws.Range("K28").Value = 1583.2354
$ws.Range("M28").Value = -1098.2354
# row 33 (chunk 2)
$ws.Range("H33").Value = 382.08694
$ws.Range("I33").Value = 383.5909
$ws.Range("K33").Value = 383.5909
$ws.Range("M33").Value = -154.5909
# row 62 (chunk 3)
$ws.Range("H62").Value = 3574.25
$ws.Range("I62").Value = 1800
$ws.Range("J62").Value = 4165.6665
$ws.Range("K62").Value = 1800
$ws.Range("L62").Value = 4165.6665
$ws.Range("M62").Value = -1176
$ws.Range("N62").Value = -5413.6665
# row 65 (chunk 4)
$ws.Range("H65").Value = 3574.25
$ws.Range("I65").Value = 1800
$ws.Range("J65").Value = 4165.6665
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 20828.3325
$ws.Range("M65").Value = -5880
$ws.Range("N65").Value = -27068.3325
# row 113 (chunk 5)
$ws.Range("H113").Value = 11000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 11000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 11000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -17508
# row 116 (chunk 6)
$ws.Range("H116").Value = 4458.143
$ws.Range("I116").Value = 1998.3334
$ws.Range("K116").Value = 1998.3334
$ws.Range("M116").Value = 1443.6666
# row 137 (chunk 7)
$ws.Range("H137").Value = 1817.75
$ws.Range("I137").Value = 1817.75
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5453.25
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -2903.25
$ws.Range("N137").ClearContents()
# row 141 (chunk 8)
$ws.Range("H141").Value = 1393.5312
$ws.Range("J141").Value = 242
$ws.Range("L141").Value = 726
$ws.Range("N141").Value = -11086

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 45 (chunk 9)
$ws.Range("H45").Value = 1728.8
$ws.Range("I45").Value = 1728.8
$ws.Range("K45").Value = 1728.8
$ws.Range("M45").Value = -1351.8
# row 61 (chunk 10)
$ws.Range("H61").Value = 4575.4
$ws.Range("I61").Value = 959.1667
$ws.Range("K61").Value = 959.1667
$ws.Range("M61").Value = -747.1667
# row 74 (chunk 11)
$ws.Range("H74").Value = 2159.8
$ws.Range("J74").Value = 1999
$ws.Range("L74").Value = 1999
$ws.Range("N74").Value = -3747
# row 77 (chunk 12)
$ws.Range("H77").Value = 2159.8
$ws.Range("J77").Value = 1999
$ws.Range("L77").Value = 9995
$ws.Range("N77").Value = -18731
# row 97 (chunk 13)
$ws.Range("H97").Value = 1318.579
$ws.Range("I97").Value = 856.05884
$ws.Range("K97").Value = 856.05884
$ws.Range("M97").Value = -360.05884
# row 110 (chunk 14)
$ws.Range("H110").Value = 7643.091
$ws.Range("I110").Value = 7771.1113
$ws.Range("J110").Value = 7067
$ws.Range("K110").Value = 7771.1113
$ws.Range("L110").Value = 7067
$ws.Range("M110").Value = -5726.1113
$ws.Range("N110").Value = -11157
# row 136 (chunk 15)
$ws.Range("H136").Value = 4575.4
$ws.Range("I136").Value = 959.1667
$ws.Range("K136").Value = 2877.5001
$ws.Range("M136").Value = -327.5001000000002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 16 (chunk 16)
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# row 29 (chunk 17)
$ws.Range("H29").Value = 6000
$ws.Range("I29").Value = 6000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -5711
$ws.Range("N29").ClearContents()
# row 80 (chunk 18)
$ws.Range("H80").Value = 4985.7334
$ws.Range("I80").Value = 610.5454999999999
$ws.Range("K80").Value = 610.5454999999999
$ws.Range("M80").Value = 387.4545000000001
# row 83 (chunk 19)
$ws.Range("H83").Value = 4985.7334
$ws.Range("I83").Value = 610.5454999999999
$ws.Range("K83").Value = 3052.7275
$ws.Range("M83").Value = 1939.2725
# row 134 (chunk 20)
$ws.Range("H134").Value = 5356.0586
$ws.Range("I134").Value = 5328.9033
$ws.Range("K134").Value = 15986.7099
$ws.Range("M134").Value = -13451.7099

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 108 (chunk 21)
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# row 134 (chunk 22)
$ws.Range("H134").Value = 2147.1738
$ws.Range("I134").Value = 2024.5714
$ws.Range("J134").Value = 3434.5
$ws.Range("K134").Value = 6073.7142
$ws.Range("L134").Value = 10303.5
$ws.Range("M134").Value = -3538.7142
$ws.Range("N134").Value = -15373.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 69 (chunk 23)
$ws.Range("H69").Value = 2999
$ws.Range("I69").Value = 2999
$ws.Range("K69").Value = 8997
$ws.Range("M69").Value = -8186
# row 72 (chunk 24)
$ws.Range("H72").Value = 2999
$ws.Range("I72").Value = 2999
$ws.Range("K72").Value = 26991
$ws.Range("M72").Value = -22935
# row 92 (chunk 25)
$ws.Range("H92").Value = 271.25
$ws.Range("I92").Value = 271.25
$ws.Range("K92").Value = 813.75
$ws.Range("M92").Value = 434.25
# row 129 (chunk 26)
$ws.Range("H129").Value = 3415.2856
$ws.Range("I129").Value = 3632.6667
$ws.Range("J129").Value = 3252.25
$ws.Range("K129").Value = 10898.0001
$ws.Range("L129").Value = 9756.75
$ws.Range("M129").Value = -5898.000100000001
$ws.Range("N129").Value = -19756.75
# row 131 (chunk 27)
$ws.Range("H131").Value = 3157.9
$ws.Range("J131").Value = 3580.5
$ws.Range("L131").Value = 10741.5
$ws.Range("N131").Value = -20821.5
# row 132 (chunk 28)
$ws.Range("H132").Value = 1799.1666
$ws.Range("I132").Value = 1899
$ws.Range("J132").Value = 1749.25
$ws.Range("K132").Value = 17091
$ws.Range("L132").Value = 15743.25
$ws.Range("M132").Value = -14561
$ws.Range("N132").Value = -20803.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 24 (chunk 29)
$ws.Range("H24").Value = 335003
$ws.Range("I24").Value = 300000
$ws.Range("J24").Value = 346670.66
$ws.Range("K24").Value = 300000
$ws.Range("L24").Value = 346670.66
$ws.Range("M24").Value = -299827
$ws.Range("N24").Value = -347016.66
# row 70 (chunk 30)
$ws.Range("H70").Value = 9178.77
$ws.Range("I70").Value = 9002.4
$ws.Range("J70").Value = 9766.666999999999
$ws.Range("K70").Value = 9002.4
$ws.Range("L70").Value = 9766.666999999999
$ws.Range("M70").Value = -8732.4
$ws.Range("N70").Value = -10306.667
# row 73 (chunk 31)
$ws.Range("H73").Value = 9178.77
$ws.Range("I73").Value = 9002.4
$ws.Range("J73").Value = 9766.666999999999
$ws.Range("K73").Value = 9002.4
$ws.Range("L73").Value = 9766.666999999999
$ws.Range("M73").Value = -8066.4
$ws.Range("N73").Value = -11638.667
# row 102 (chunk 32)
$ws.Range("H102").Value = 1077.75
$ws.Range("I102").Value = 1270.3334
$ws.Range("K102").Value = 1270.3334
$ws.Range("M102").Value = 351.6666
# row 113 (chunk 33)
$ws.Range("H113").Value = 3510
$ws.Range("I113").Value = 3010
$ws.Range("J113").Value = 4010
$ws.Range("K113").Value = 3010
$ws.Range("L113").Value = 4010
$ws.Range("M113").Value = -840
$ws.Range("N113").Value = -8350
# row 126 (chunk 34)
$ws.Range("H126").Value = 3924.3333
$ws.Range("I126").Value = 3924.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11772.9999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9302.999899999999
$ws.Range("N126").ClearContents()
# row 132 (chunk 35)
$ws.Range("H132").Value = 2654.0715
$ws.Range("I132").Value = 2344.1667
$ws.Range("K132").Value = 7032.500100000001
$ws.Range("M132").Value = -4502.500100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 76 (chunk 36)
$ws.Range("H76").Value = 12288
$ws.Range("J76").Value = 12288
$ws.Range("L76").Value = 12288
$ws.Range("N76").Value = -12964
# row 79 (chunk 37)
$ws.Range("H79").Value = 12288
$ws.Range("J79").Value = 12288
$ws.Range("L79").Value = 12288
$ws.Range("N79").Value = -14628
# row 100 (chunk 38)
$ws.Range("H100").Value = 2034.3334
$ws.Range("I100").Value = 2001.5
$ws.Range("K100").Value = 2001.5
$ws.Range("M100").Value = -1460.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 20 (chunk 39)
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9760
# row 75 (chunk 40)
$ws.Range("H75").Value = 44999.332
$ws.Range("I75").Value = 35000
$ws.Range("K75").Value = 35000
$ws.Range("M75").Value = -34064
# row 78 (chunk 41)
$ws.Range("H78").Value = 44999.332
$ws.Range("I78").Value = 35000
$ws.Range("K78").Value = 105000
$ws.Range("M78").Value = -100320
# row 132 (chunk 42)
$ws.Range("H132").Value = 3898.5557
$ws.Range("J132").Value = 5895
$ws.Range("L132").Value = 17685
$ws.Range("N132").Value = -22745
